$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "tabla de datos dispersos"

# Preserve the outline properties (sheetPr/outlinePr) on the new data sheet
$ws2.Outline.SummaryBelow = $true
$ws2.Outline.SummaryRight = $true

# Copy all used-range content (values, formats, merges) from ws1 to ws2
$ws1.UsedRange.Copy($ws2.Range("A1"))

# Clear ws1 entirely (contents, formats, and merged cells)
$ws1.Cells.UnMerge()
$ws1.Cells.Clear()

# Rename ws1
$ws1.Name = "Sheet"

Write-Output $wb.Worksheets.Count
Write-Output $wb.Worksheets.Item(1).Name
Write-Output $wb.Worksheets.Item(2).Name
